# Auto-generated edit script: refresh market-price derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job
# Leve-profit tables, mirroring a scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 7475
$ws.Range("I74").Value = 7933.3335
$ws.Range("J74").Value = 6100
$ws.Range("K74").Value = 7933.3335
$ws.Range("L74").Value = 6100
$ws.Range("M74").Value = -6997.3335
$ws.Range("N74").Value = -7972

$ws.Range("H77").Value = 7475
$ws.Range("I77").Value = 7933.3335
$ws.Range("J77").Value = 6100
$ws.Range("K77").Value = 39666.6675
$ws.Range("L77").Value = 30500
$ws.Range("M77").Value = -34986.6675
$ws.Range("N77").Value = -39860

$ws.Range("H99").Value = 254.33333
$ws.Range("I99").Value = 300
$ws.Range("J99").Value = 231.5
$ws.Range("K99").Value = 900
$ws.Range("L99").Value = 694.5
$ws.Range("M99").Value = 598
$ws.Range("N99").Value = -3690.5

$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()

$ws.Range("H129").Value = 1869.5
$ws.Range("J129").Value = 2244.5
$ws.Range("L129").Value = 6733.5
$ws.Range("N129").Value = -16733.5

$ws.Range("H137").Value = 1944.3334
$ws.Range("I137").Value = 1974.875
$ws.Range("K137").Value = 5924.625
$ws.Range("M137").Value = -3374.625

$ws.Range("H138").Value = 2285.818
$ws.Range("I138").Value = 1372.4286
$ws.Range("J138").Value = 3884.25
$ws.Range("K138").Value = 4117.2858
$ws.Range("L138").Value = 11652.75
$ws.Range("M138").Value = 1022.7142
$ws.Range("N138").Value = -21932.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2535.2856
$ws.Range("I61").Value = 1499.3334
$ws.Range("K61").Value = 1499.3334
$ws.Range("M61").Value = -1287.3334

$ws.Range("H63").Value = 11501.333
$ws.Range("I63").Value = 7249.25
$ws.Range("J63").Value = 20005.5
$ws.Range("K63").Value = 7249.25
$ws.Range("L63").Value = 20005.5
$ws.Range("M63").Value = -6563.25
$ws.Range("N63").Value = -21377.5

$ws.Range("H66").Value = 11501.333
$ws.Range("I66").Value = 7249.25
$ws.Range("J66").Value = 20005.5
$ws.Range("K66").Value = 36246.25
$ws.Range("L66").Value = 100027.5
$ws.Range("M66").Value = -32814.25
$ws.Range("N66").Value = -106891.5

$ws.Range("H74").Value = 837.35297
$ws.Range("I74").Value = 820.0625
$ws.Range("J74").Value = 1114
$ws.Range("K74").Value = 820.0625
$ws.Range("L74").Value = 1114
$ws.Range("M74").Value = 53.9375
$ws.Range("N74").Value = -2862

$ws.Range("H77").Value = 837.35297
$ws.Range("I77").Value = 820.0625
$ws.Range("J77").Value = 1114
$ws.Range("K77").Value = 4100.3125
$ws.Range("L77").Value = 5570
$ws.Range("M77").Value = 267.6875
$ws.Range("N77").Value = -14306

$ws.Range("H97").Value = 2437.2
$ws.Range("I97").Value = 2596.889
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 2596.889
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -2100.889
$ws.Range("N97").Value = -1992

$ws.Range("H102").Value = 2560.5
$ws.Range("I102").Value = 2899
$ws.Range("J102").Value = 2222
$ws.Range("K102").Value = 2899
$ws.Range("L102").Value = 2222
$ws.Range("M102").Value = -1277
$ws.Range("N102").Value = -5466

$ws.Range("H132").Value = 2537.8948
$ws.Range("I132").Value = 1401.6923
$ws.Range("K132").Value = 4205.0769
$ws.Range("M132").Value = -1675.0769

$ws.Range("H136").Value = 2535.2856
$ws.Range("I136").Value = 1499.3334
$ws.Range("K136").Value = 4498.0002
$ws.Range("M136").Value = -1948.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 945.36365
$ws.Range("I94").Value = 939.9
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 939.9
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -488.9
$ws.Range("N94").Value = -1902

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").ClearContents()

$ws.Range("H134").Value = 3895.652
$ws.Range("I134").Value = 4971.846
$ws.Range("J134").Value = 2496.6
$ws.Range("K134").Value = 14915.538
$ws.Range("L134").Value = 7489.799999999999
$ws.Range("M134").Value = -12380.538
$ws.Range("N134").Value = -12559.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 724.5
$ws.Range("I6").Value = 450
$ws.Range("J6").Value = 999
$ws.Range("K6").Value = 450
$ws.Range("L6").Value = 999
$ws.Range("M6").Value = -337
$ws.Range("N6").Value = -1225

$ws.Range("H7").Value = 218.58333
$ws.Range("I7").Value = 152.4
$ws.Range("K7").Value = 152.4
$ws.Range("M7").Value = -39.40000000000001

$ws.Range("H31").Value = 66669588
$ws.Range("I31").Value = 125001096
$ws.Range("K31").Value = 125001096
$ws.Range("M31").Value = -125000801

$ws.Range("H34").Value = 66669588
$ws.Range("I34").Value = 125001096
$ws.Range("K34").Value = 125001096
$ws.Range("M34").Value = -125000894

$ws.Range("H86").Value = 3699.75
$ws.Range("I86").Value = 3699.75
$ws.Range("K86").Value = 3699.75
$ws.Range("M86").Value = -2576.75

$ws.Range("H89").Value = 3699.75
$ws.Range("I89").Value = 3699.75
$ws.Range("K89").Value = 18498.75
$ws.Range("M89").Value = -12882.75

$ws.Range("H105").Value = 3625
$ws.Range("I105").Value = 3350
$ws.Range("K105").Value = 3350
$ws.Range("M105").Value = -1603

$ws.Range("H132").Value = 3369.2222
$ws.Range("I132").Value = 3265.2
$ws.Range("J132").Value = 3499.25
$ws.Range("K132").Value = 9795.599999999999
$ws.Range("L132").Value = 10497.75
$ws.Range("M132").Value = -7265.599999999999
$ws.Range("N132").Value = -15557.75

$ws.Range("H134").Value = 2961.4285
$ws.Range("I134").Value = 2961.4285
$ws.Range("K134").Value = 8884.2855
$ws.Range("M134").Value = -6349.2855

$ws.Range("H138").Value = 68995.8
$ws.Range("J138").Value = 68995.8
$ws.Range("L138").Value = 68995.8
$ws.Range("N138").Value = -79275.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 55.8
$ws.Range("I23").Value = 55.8
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 167.4
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 67.60000000000002
$ws.Range("N23").ClearContents()

$ws.Range("H132").Value = 3250
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3250
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 29250
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -34310

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2974.75
$ws.Range("J80").Value = 2966.3333
$ws.Range("L80").Value = 2966.3333
$ws.Range("N80").Value = -4962.3333

$ws.Range("H83").Value = 2974.75
$ws.Range("J83").Value = 2966.3333
$ws.Range("L83").Value = 14831.6665
$ws.Range("N83").Value = -24815.6665

$ws.Range("H97").Value = 322
$ws.Range("I97").Value = 306.7
$ws.Range("J97").Value = 373
$ws.Range("K97").Value = 306.7
$ws.Range("L97").Value = 373
$ws.Range("M97").Value = 189.3
$ws.Range("N97").Value = -1365

$ws.Range("H122").Value = 2694.4
$ws.Range("I122").Value = 2784
$ws.Range("K122").Value = 8352
$ws.Range("M122").Value = -5902

$ws.Range("H132").Value = 2941.4167
$ws.Range("I132").Value = 2255.5557
$ws.Range("K132").Value = 6766.6671
$ws.Range("M132").Value = -4236.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1862.75
$ws.Range("I93").Value = 1385.5714
$ws.Range("J93").Value = 2530.8
$ws.Range("K93").Value = 1385.5714
$ws.Range("L93").Value = 2530.8
$ws.Range("M93").Value = -137.5714
$ws.Range("N93").Value = -5026.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1000
$ws.Range("I96").Value = 1000
$ws.Range("K96").Value = 1000
$ws.Range("M96").Value = 373

$ws.Range("H100").Value = 3000
$ws.Range("I100").Value = 3000
$ws.Range("K100").Value = 6000
$ws.Range("M100").Value = -5459

$ws.Range("H122").Value = 55557110
$ws.Range("I122").Value = 1785.7142
$ws.Range("J122").Value = 250000750
$ws.Range("K122").Value = 5357.142599999999
$ws.Range("L122").Value = 750002250
$ws.Range("M122").Value = -2907.142599999999
$ws.Range("N122").Value = -750007150

$ws.Range("H126").Value = 1464
$ws.Range("J126").Value = 832.6667
$ws.Range("L126").Value = 2498.0001
$ws.Range("N126").Value = -7438.0001

$ws.Range("H132").Value = 2473.2083
$ws.Range("J132").Value = 4389.7
$ws.Range("L132").Value = 13169.1
$ws.Range("N132").Value = -18229.1
